$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - updated values from the new MATLAB run
$ws.Range("A2").Value = 0.63867721557144308
$ws.Range("B2").Value = 1.9416657838641791
$ws.Range("C2").Value = 0.0031596877710231564
$ws.Range("D2").Value = 0.067205378687909836
$ws.Range("E2").Value = 0.68802191628728926
$ws.Range("F2").Value = 0.13848281731183243
$ws.Range("G2").Value = 49.600000000000001
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 1.181359437926748
$ws.Range("J2").Value = 0.096194337846943451
$ws.Range("K2").Value = 0.48580736052526474
$ws.Range("L2").Value = -0.24572604525856168
$ws.Range("M2").Value = -0.95165831459267314

# Row 3 - new row
$ws.Range("A3").Value = 0.56178793934858784
$ws.Range("B3").Value = 2.1383016272951534
$ws.Range("C3").Value = 0.0017549989067366933
$ws.Range("D3").Value = 0.48566544174909315
$ws.Range("E3").Value = 0.6874499523968628
$ws.Range("F3").Value = 0.13848281731183243
$ws.Range("G3").Value = 49.600000000000001
$ws.Range("H3").Value = 4
$ws.Range("I3").Value = 0.81288935844146215
$ws.Range("J3").Value = -0.055605437034776356
$ws.Range("K3").Value = 0.27628327593547214
$ws.Range("L3").Value = 0.11330974531366156
$ws.Range("M3").Value = -0.95246903675981165

# Row 4 - new row
$ws.Range("A4").Value = 0.56154735255209332
$ws.Range("B4").Value = 2.1388885309813892
$ws.Range("C4").Value = 0.011132807693144287
$ws.Range("D4").Value = 0.48566544174909315
$ws.Range("E4").Value = 0.6874499523968628
$ws.Range("F4").Value = 0.13848281731183243
$ws.Range("G4").Value = 49.600000000000001
$ws.Range("H4").Value = 10
$ws.Range("I4").Value = 0.81288935844146215
$ws.Range("J4").Value = -0.04641284912580082
$ws.Range("K4").Value = 0.2857649451163331
$ws.Range("L4").Value = 0.12051995808458571
$ws.Range("M4").Value = -0.94461671108638257
